$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"

$ws.Range("B4").Select()
